{"js": "// Find the two target paragraphs by matching their (pre-edit) text so the\n// script is robust even if surrounding content shifts a little.\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\nconst siftText =\n  \"Sift through 250 um sieve, anything that goes through gets put in a jar. Anything that remains in the sieve needs to be re-ground.\";\nconst pomText =\n  \"This is to make sure the particle size is representative of particulate organic matter (POM) that microbes would experience in their natural environment.\";\n\nlet siftPara = null;\nlet pomPara = null;\nfor (let i = 0; i < paras.items.length; i++) {\n  const t = paras.items[i].text;\n  if (t === siftText) siftPara = paras.items[i];\n  if (t === pomText) pomPara = paras.items[i];\n}\n\nif (!siftPara || !pomPara) {\n  throw new Error(\"Could not locate target paragraphs\");\n}\n\n// 1. Shorten the \"Sift through...\" sentence: drop the trailing clause that\n//    is about to become its own bullet.\nsiftPara.insertText(\n  \"Sift through 250 um sieve, anything that goes through gets put in a jar.\",\n  \"Replace\"\n);\n\n// 2. The former \"particle size\" sub-bullet becomes the re-grind instruction\n//    (it already sits at the right list/level, so only its text changes).\npomPara.insertText(\n  \"Anything that remains in the sieve needs to be re-ground.\",\n  \"Replace\"\n);\n\n// 3. Add a brand-new sub-bullet right after it carrying the reworded\n//    \"sieving\" explanation, matching the same list (numId) and level.\nconst newPara = pomPara.insertParagraph(\n  \"The sieving is to make sure the particle size is representative of particulate organic matter (POM) that microbes would experience in their natural environment.\",\n  \"After\"\n);\nnewPara.attachToList(1003, 1);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the two paragraphs we need by their current text so the script\n# does not depend on a hard-coded paragraph index.\n$siftText = \"Sift through 250 um sieve, anything that goes through gets put in a jar. Anything that remains in the sieve needs to be re-ground.\"\n$pomText  = \"This is to make sure the particle size is representative of particulate organic matter (POM) that microbes would experience in their natural environment.\"\n\n$siftPara = $null\n$pomPara  = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd(\"`r\", \"`n\", \"`v\")\n    if ($t -eq $siftText) { $siftPara = $p }\n    if ($t -eq $pomText)  { $pomPara  = $p }\n}\n\n# 1. Shorten the \"Sift through...\" sentence: the trailing clause is about\n#    to become its own bullet below.\n$siftPara.Range.Text = \"Sift through 250 um sieve, anything that goes through gets put in a jar.\"\n\n# 2. The former \"particle size\" sub-bullet becomes the re-grind\n#    instruction (it already sits at the right list/level, so only its\n#    text changes).\n$pomPara.Range.Text = \"Anything that remains in the sieve needs to be re-ground.\"\n\n# 3. Add a brand-new sub-bullet right after it carrying the reworded\n#    \"sieving\" explanation; InsertParagraphAfter() clones the paragraph\n#    (and list/numbering) formatting of $pomPara automatically.\n$newRange = $pomPara.Range.InsertParagraphAfter()\n$newPara = $pomPara.Next()\n$newPara.Range.Text = \"The sieving is to make sure the particle size is representative of particulate organic matter (POM) that microbes would experience in their natural environment.\"\n\n$d.Save()\n"}
